# RxSwiftBasics3 - Day3 slide: add "Scan Operator" to the Day 3 bullet,
# merge the leading space into the Day 4 bullet's run, and refresh the
# RxDataSources run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Locate the body placeholder (ppPlaceholderBody = 2) that lists the Day 1-5
# agenda, rather than assuming a fixed shape index.
$shp = $null
$phs = $s.Shapes.Placeholders
for ($i = 1; $i -le $phs.Count; $i++) {
    $candidate = $phs.Item($i)
    if ($candidate.PlaceholderFormat.Type -eq 2) {
        $shp = $candidate
    }
}
if ($shp -eq $null) {
    $shp = $s.Shapes.Item(4)
}

$tf = $shp.TextFrame
$tr = $tf.TextRange

# ---------------------------------------------------------------------
# 1) Day 3 bullet: "...hide 'Loading' )" -> "...hide 'Loading' ), Scan Operator"
#    Split so the trailing ")" becomes "), Scan Operator", then split the
#    closing quote + space into its own run.
# ---------------------------------------------------------------------
$full = $tr.Text
$loadingIdx = $full.IndexOf("Loading")
$parenPos = $loadingIdx + 9              # 0-based index of the closing ")"
$parenRun = $tr.Characters($parenPos + 1, 1)
$parenRun.Text = "), Scan Operator"

$full = $tr.Text
$loadingIdx = $full.IndexOf("Loading")
$quotePos = $loadingIdx + 7              # 0-based index of the closing quote
$quoteRun = $tr.Characters($quotePos + 1, 2)
$quoteRun.Text = [char]0x2019 + " "

# ---------------------------------------------------------------------
# 2) Day 4 bullet: merge the standalone " " run into the following run so
#    the line reads "... Adding a Reactive Extension to Custom UI Element, "
#    as a single run.
# ---------------------------------------------------------------------
$full = $tr.Text
$addingIdx = $full.IndexOf("Adding a Reactive")
$spanStart = $addingIdx - 1
$addingText = "Adding a Reactive Extension to Custom UI Element, "
$spanLen = 1 + $addingText.Length
$mergedRun = $tr.Characters($spanStart + 1, $spanLen)
$mergedRun.Text = " " + $addingText

# ---------------------------------------------------------------------
# 3) RxDataSources run: touch it so it is refreshed/re-saved.
# ---------------------------------------------------------------------
$full = $tr.Text
$rxIdx = $full.IndexOf("RxDataSources")
$rxRun = $tr.Characters($rxIdx + 1, "RxDataSources".Length)
$rxRun.Text = "RxDataSources"
